$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values: L1 "Lunch" -> "Fomal", O1 "Music" -> "Non-Music"
$ws.Range("L1").Value = "Fomal"
$ws.Range("O1").Value = "Non-Music"

# Column widths (E and O)
$ws.Columns.Item(5).ColumnWidth = 13.5
$ws.Columns.Item(15).ColumnWidth = 10.333333333333334

# Update selection to O1
$ws.Range("O1").Select()
